# Updated data model upgrade description and implemented consistent
# abbreviations in object names (DB / PL-SQL naming conventions).
#
# The "View Comments" sheet documents column comments for two source
# views, CCD_DATA_SETS_V (rows 496-509) and DS_PIR_SCOR_V (rows 510-550).
# Both views were renamed/merged into CCD_DATA_SETS_INPORT_V, so column A
# for rows 496-550 needs to point at that new name. Column D is a
# CONCATENATE() formula that already references column A, so it
# recalculates automatically once A is updated.

$wb = $excel.ActiveWorkbook

# --- "View Query Builder" sheet: selection moves, no data changes ------
$wsBuilder = $wb.Worksheets.Item("View Query Builder")
$wsBuilder.Activate() | Out-Null
$wsBuilder.Range("C1260:C1301").Select() | Out-Null

# --- "View Comments" sheet: rename source view in column A -------------
$wsComments = $wb.Worksheets.Item("View Comments")
$wsComments.Activate() | Out-Null
$wsComments.Range("A496:A550").Value = "CCD_DATA_SETS_INPORT_V"

# Restore the selection/active cell shown in the sheet after the edit.
$wsComments.Range("B290").Select() | Out-Null
